$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the login-link cell (A4) to point at the new staging URL.
$ws.Range("A4").Value = "https://stg.oxs.co.il/"

# Replicate the author's saved cursor position (B12) from the commit.
$null = $ws.Range("B12").Select()
